# Update odds values on the active worksheet (Sheet1) to reflect the
# latest FlashScore odds for 2024-10-16 matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Fortaleza vs Atletico-MG)
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65

# Row 8 (Correcaminos vs Cancun)
$ws.Range("G8").Value = 3.15
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 2.18
$ws.Range("K8").Value = 2.12
$ws.Range("L8").Value = 2.72
$ws.Range("O8").Value = 1.31
$ws.Range("P8").Value = 2.9
$ws.Range("U8").Value = 1.72
$ws.Range("V8").Value = 1.88
$ws.Range("W8").Value = 9.5
$ws.Range("X8").Value = 16.5
$ws.Range("AC8").Value = 9.25
$ws.Range("AD8").Value = 6.3
$ws.Range("AE8").Value = 14
$ws.Range("AF8").Value = 65
$ws.Range("AG8").Value = 500
$ws.Range("AH8").Value = 7.3
$ws.Range("AJ8").Value = 8.75
$ws.Range("AK8").Value = 21
$ws.Range("AL8").Value = 18
$ws.Range("AN8").Value = 5.1
$ws.Range("AP8").Value = 22
$ws.Range("AR8").Value = 100
$ws.Range("AS8").Value = 250
$ws.Range("AU8").Value = 6.7
$ws.Range("AV8").Value = 55
$ws.Range("AW8").Value = 4.1
$ws.Range("AX8").Value = 11
$ws.Range("AZ8").Value = 45
$ws.Range("BB8").Value = 200
